$d = $word.ActiveDocument

# The "Download chart data and load to database / view" bullet currently sits
# right after "Download company data and load to database / view". It needs
# to move up so it sits right before "Watchlist view / form" instead.

$chartText = "Download chart data and load to database / view"

# Step 1: insert a new paragraph (matching the list formatting of its new
# neighbor) right before "Watchlist view / form", and give it the moved text.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Watchlist view / form`r") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newPara.Range.Text = $chartText
        break
    }
}

# Step 2: remove the original "Download chart data..." paragraph that used
# to follow "Download company data and load to database / view".
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Download company data and load to database / view`r") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -eq "$chartText`r") {
            $nextPara.Range.Delete()
        }
        break
    }
}
